$d = $word.ActiveDocument

# Locate the paragraph that contains the "standalone.sh" command line; the
# new content is inserted right after it (and before the two trailing blank
# paragraphs already present at the end of the document).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "standalone.sh*") {
        $target = $p
    }
}
if ($target -eq $null) {
    throw "Could not find the 'standalone.sh' paragraph"
}

$insertPos = $target.Range.End
$r = $d.Range($insertPos, $insertPos)

$xml = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t>Contenedor</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">Para que </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>JBoss</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> reconozca los caracteres UTF para la subida de documentos (Quizá en más aspectos) es necesario cambiar una opción en el fichero standalone.xml</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">Para solucionarlo hay que editar el </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
        </w:rPr>
        <w:t>standalone.xml</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> y sustituir la línea </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>&lt;</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>servlet-container</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>name</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>="default”&gt;</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">  por </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>&lt;</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>servlet-container</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>name</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>="default" default-</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>encoding</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>="UTF-8"&gt;</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>

'@

$r.InsertXML($xml)

Write-Output "Paragraph count after insert:"
Write-Output $d.Paragraphs.Count
